# On PictureInstructionsVera: Adding files which were renamed so that CC65
# doesn't delete them to project Allocation Work.
#
# This adds a new "Vera" worksheet (VRAM budget table) at the end of the
# workbook and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- Create the new "Vera" worksheet as the last tab ---------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Vera"

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "Purpose"
$ws.Range("B1").Value = "Start"
$ws.Range("C1").Value = "Size"
$ws.Range("D1").Value = "Vera Address"
$ws.Range("E1").Value = "Ends"
$ws.Range("A1:E1").Font.Bold = $true

# --- Row 2 : Background Images ---------------------------------------------
$ws.Range("A2").Value = "Background Images"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 38400
$ws.Range("D2").Formula = '=DEC2HEX(HEX2DEC("00000"))'
$ws.Range("E2").Formula = '=DEC2HEX(HEX2DEC(D2)+C2-1)'

# --- Row 3 : Priority Screen ------------------------------------------------
$ws.Range("A3").Value = "Priority Screen"
$ws.Range("B3").Formula = '=C2'
$ws.Range("C3").Value = 13360
$ws.Range("D3").Formula = '=DEC2HEX(B2 + C2)'
$ws.Range("E3").Formula = '=DEC2HEX(HEX2DEC(D3)+C3-1)'

# --- Row 4 : Vacant ----------------------------------------------------------
$ws.Range("A4").Value = "Vacant"
$ws.Range("B4").Formula = '=B3+C3'
$ws.Range("C4").Value = 1488
$ws.Range("D4").Formula = '=DEC2HEX(B3 + C3)'
$ws.Range("E4").Formula = '=DEC2HEX(HEX2DEC(D4)+C4-1)'

# --- Row 5 : Tilebase --------------------------------------------------------
$ws.Range("A5").Value = "Tilebase"
$ws.Range("B5").Formula = '=B4+C4'
$ws.Range("C5").Value = 2560
$ws.Range("D5").Formula = '=DEC2HEX(B4 + C4)'
$ws.Range("E5").Formula = '=DEC2HEX(HEX2DEC(D5)+C5-1)'

# --- Row 6 : Map Base --------------------------------------------------------
$ws.Range("A6").Value = "Map Base"
$ws.Range("B6").Formula = '=B5+C5'
$ws.Range("C6").Value = 2561
$ws.Range("D6").Formula = '=DEC2HEX(B5 + C5)'
$ws.Range("E6").Formula = '=DEC2HEX(HEX2DEC(D6)+C6-1)'

# --- Row 7 : Vacant -----------------------------------------------------------
$ws.Range("A7").Value = "Vacant"
$ws.Range("B7").Formula = '=B6+C6'
$ws.Range("C7").Value = 1023
$ws.Range("D7").Formula = '=DEC2HEX(B6 + C6)'
$ws.Range("E7").Formula = '=DEC2HEX(HEX2DEC(D7)+C7-1)'

# --- Totals (entered before rows 8/9 so new shared-string indices line up
#     with the authored order: Total Used/Available/Remaining precede
#     Sprites/Volatile Buffer in xl/sharedStrings.xml) ------------------------
$ws.Range("B12").Value = "Total Used"
$ws.Range("C12").Formula = '=SUM(C2:C9)'

$ws.Range("B13").Value = "Total Available"
$ws.Range("C13").Value = 129471

$ws.Range("B14").Value = "Remaining"
$ws.Range("C14").Formula = '=C13-C12'

# --- Row 8 : Sprites -----------------------------------------------------------
$ws.Range("A8").Value = "Sprites"
$ws.Range("B8").Formula = '=B7+C7'
$ws.Range("C8").Value = 60079
$ws.Range("D8").Formula = '=DEC2HEX(B7 + C7)'
$ws.Range("E8").Formula = '=DEC2HEX(HEX2DEC(D8)+C8-1)'

# --- Row 9 : Volatile Buffer -----------------------------------------------------
$ws.Range("A9").Value = "Volatile Buffer"
$ws.Range("B9").Formula = '=B8+C8'
$ws.Range("C9").Value = 10000
$ws.Range("D9").Formula = '=DEC2HEX(B8 + C8)'
$ws.Range("E9").Formula = '=DEC2HEX(HEX2DEC(D9)+C9-1)'

# --- Column widths (cosmetic, best effort) --------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.45
$ws.Range("D1:E1").ColumnWidth = 15.82

# --- Selection / active cell, matches authored state ----------------------------
[void]$ws.Range("D9").Select()

# --- Make the new sheet the active tab ------------------------------------------
$ws.Activate()
